$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_LT"

$values = @(88.036809815950917, 88.343558282208591, 89.033742331288352, 91.027607361963192, 91.334355828220865, 91.334355828220865, 78.757668711656436, 80.981595092024534, 81.058282208588963, 78.911042944785279, 81.211656441717793, 83.819018404907979, 79.907975460122699, 79.831288343558285, 82.208588957055213, 80.828220858895705, 83.588957055214721, 87.269938650306742, 89.340490797546011, 89.723926380368098, 89.417177914110425, 86.963190184049083, 89.953987730061343, 89.800613496932513, 84.662576687116569, 84.432515337423311, 84.815950920245399, 81.901840490797554, 81.058282208588963, 81.978527607361968, 88.113496932515332, 93.865030674846622, 92.407975460122699, 86.349693251533751, 93.558282208588963, 83.742331288343564, 90.030674846625772, 86.042944785276077, 86.886503067484668, 81.441717791411037, 81.595092024539866, 80.138036809815944, 80.904907975460134, 85.122699386503058, 84.355828220858896, 80.981595092024534, 81.748466257668724, 81.058282208588963)
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
